# Update cfb_weather workbook data (FBS and Other sheets) with refreshed weather/odds data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# ---- FBS sheet ----
# Row 2
$ws1.Cells.Item(2, 1).Value = 'Troy @ James Madison'
$ws1.Cells.Item(2, 2).Value = 'FRI 12/05'
$ws1.Cells.Item(2, 3).Value = '07:00 PM'
$ws1.Cells.Item(2, 4).Value = 'High'
$ws1.Cells.Item(2, 5).Value = 'N-S'
$ws1.Cells.Item(2, 6).Value = 'Med'
$ws1.Cells.Item(2, 7).Value = 'E'
$ws1.Cells.Item(2, 8).Value = 244.5283813
$ws1.Cells.Item(2, 9).Value = 54.52
$ws1.Cells.Item(2, 10).Value = 65.95999999999999
$ws1.Cells.Item(2, 11).Value = 6.2
$ws1.Cells.Item(2, 12).Value = 1975
$ws1.Cells.Item(2, 13).Value = 'NNE'
$ws1.Cells.Item(2, 14).Value = 'NNE'
$ws1.Cells.Item(2, 15).Value = 28.94
$ws1.Cells.Item(2, 16).Value = 5.5
$ws1.Cells.Item(2, 17).Value = 'NNE'
$ws1.Cells.Item(2, 18).Value = 0
$ws1.Cells.Item(2, 19).Value = -0.13
$ws1.Cells.Item(2, 20).Value = -0.38
$ws1.Cells.Item(2, 21).Value = -0.7
$ws1.Cells.Item(2, 22).Value = '38.4352919, -78.8729349'
$ws1.Cells.Item(2, 23).Value = 47.5
$ws1.Cells.Item(2, 24).Value = -105
$ws1.Cells.Item(2, 25).Value = 47.5
$ws1.Cells.Item(2, 26).Value = -115
$ws1.Cells.Item(2, 27).Value = -22
$ws1.Cells.Item(2, 28).Value = -22.5
$ws1.Cells.Item(2, 31).Value = 0
$ws1.Cells.Item(2, 32).Value = 0.5
$ws1.Cells.Item(2, 37).Value = '2025-12-03T05:15:42.459576'

# Row 3
$ws1.Cells.Item(3, 1).Value = 'Kennesaw State @ Jacksonville State'
$ws1.Cells.Item(3, 2).Value = 'FRI 12/05'
$ws1.Cells.Item(3, 3).Value = '06:00 PM'
$ws1.Cells.Item(3, 4).Value = 'Low'
$ws1.Cells.Item(3, 5).Value = 'E-W'
$ws1.Cells.Item(3, 6).Value = 'High'
$ws1.Cells.Item(3, 7).Value = 'N'
$ws1.Cells.Item(3, 8).Value = -98.89108280000002
$ws1.Cells.Item(3, 9).Value = 63.15
$ws1.Cells.Item(3, 10).Value = 61.32
$ws1.Cells.Item(3, 11).Value = 4.8
$ws1.Cells.Item(3, 12).Value = 1947
$ws1.Cells.Item(3, 13).Value = 'WNW'
$ws1.Cells.Item(3, 14).Value = 'WNW'
$ws1.Cells.Item(3, 15).Value = 45.2
$ws1.Cells.Item(3, 16).Value = 1.7
$ws1.Cells.Item(3, 17).Value = 'WNW'
$ws1.Cells.Item(3, 18).Value = 0.3
$ws1.Cells.Item(3, 19).Value = 0
$ws1.Cells.Item(3, 20).Value = 0
$ws1.Cells.Item(3, 21).Value = -3.1
$ws1.Cells.Item(3, 22).Value = '33.8201052, -85.76647'
$ws1.Cells.Item(3, 23).Value = 58.5
$ws1.Cells.Item(3, 24).Value = -105
$ws1.Cells.Item(3, 25).Value = 59.5
$ws1.Cells.Item(3, 26).Value = -110
$ws1.Cells.Item(3, 27).Value = -1.5
$ws1.Cells.Item(3, 28).Value = 2.5
$ws1.Cells.Item(3, 31).Value = 0.0170940170940171
$ws1.Cells.Item(3, 32).Value = -4
$ws1.Cells.Item(3, 37).Value = '2025-12-03T05:15:42.459576'

# Row 4
$ws1.Cells.Item(4, 1).Value = 'North Texas @ Tulane'
$ws1.Cells.Item(4, 2).Value = 'FRI 12/05'
$ws1.Cells.Item(4, 3).Value = '07:00 PM'
$ws1.Cells.Item(4, 4).Value = 'High'
$ws1.Cells.Item(4, 5).Value = 'NE-SW'
$ws1.Cells.Item(4, 6).Value = 'Med'
$ws1.Cells.Item(4, 7).Value = 'NW'
$ws1.Cells.Item(4, 8).Value = -198.6871948
$ws1.Cells.Item(4, 9).Value = 70.11
$ws1.Cells.Item(4, 10).Value = 66.45
$ws1.Cells.Item(4, 11).Value = 10.6
$ws1.Cells.Item(4, 12).Value = 2014
$ws1.Cells.Item(4, 13).Value = 'SSW'
$ws1.Cells.Item(4, 14).Value = 'S'
$ws1.Cells.Item(4, 15).Value = 57.2
$ws1.Cells.Item(4, 16).Value = 6.6
$ws1.Cells.Item(4, 17).Value = 'SSW'
$ws1.Cells.Item(4, 18).Value = 0
$ws1.Cells.Item(4, 19).Value = 0
$ws1.Cells.Item(4, 20).Value = 0
$ws1.Cells.Item(4, 21).Value = -4
$ws1.Cells.Item(4, 22).Value = '29.944616, -90.116692'
$ws1.Cells.Item(4, 23).Value = 67.5
$ws1.Cells.Item(4, 24).Value = -115
$ws1.Cells.Item(4, 25).Value = 66.5
$ws1.Cells.Item(4, 26).Value = -114
$ws1.Cells.Item(4, 27).Value = 2.5
$ws1.Cells.Item(4, 28).Value = 2.5
$ws1.Cells.Item(4, 31).Value = -0.01481481481481482
$ws1.Cells.Item(4, 32).Value = 0
$ws1.Cells.Item(4, 37).Value = '2025-12-03T05:15:42.459576'

# Row 5
$ws1.Cells.Item(5, 1).Value = 'UNLV @ Boise State'
$ws1.Cells.Item(5, 2).Value = 'FRI 12/05'
$ws1.Cells.Item(5, 3).Value = '06:00 PM'
$ws1.Cells.Item(5, 4).Value = 'High'
$ws1.Cells.Item(5, 5).Value = 'N-S'
$ws1.Cells.Item(5, 6).Value = 'Med'
$ws1.Cells.Item(5, 7).Value = 'E'
$ws1.Cells.Item(5, 9).Value = 53.65
$ws1.Cells.Item(5, 10).Value = 70.04000000000001
$ws1.Cells.Item(5, 11).Value = 6.8
$ws1.Cells.Item(5, 12).Value = 1970
$ws1.Cells.Item(5, 13).Value = 'W'
$ws1.Cells.Item(5, 14).Value = 'S'
$ws1.Cells.Item(5, 15).Value = 45.44
$ws1.Cells.Item(5, 16).Value = 6.2
$ws1.Cells.Item(5, 17).Value = 'ESE'
$ws1.Cells.Item(5, 18).Value = 0
$ws1.Cells.Item(5, 19).Value = 0
$ws1.Cells.Item(5, 20).Value = 0
$ws1.Cells.Item(5, 21).Value = -0.6
$ws1.Cells.Item(5, 22).Value = '43.6028839, -116.1958882'
$ws1.Cells.Item(5, 23).Value = 57.5
$ws1.Cells.Item(5, 24).Value = -110
$ws1.Cells.Item(5, 25).Value = 58.5
$ws1.Cells.Item(5, 26).Value = -112
$ws1.Cells.Item(5, 27).Value = -3.5
$ws1.Cells.Item(5, 28).Value = -4
$ws1.Cells.Item(5, 31).Value = 0.01739130434782609
$ws1.Cells.Item(5, 32).Value = 0.5
$ws1.Cells.Item(5, 37).Value = '2025-12-03T05:15:42.459576'
$ws1.Cells.Item(5, 8).ClearContents()

# Row 6
$ws1.Cells.Item(6, 1).Value = 'Miami (OH) @ Western Michigan'
$ws1.Cells.Item(6, 2).Value = 'SAT 12/06'
$ws1.Cells.Item(6, 3).Value = '12:00 PM'
$ws1.Cells.Item(6, 4).Value = 'Mid'
$ws1.Cells.Item(6, 5).Value = 'NE-SW'
$ws1.Cells.Item(6, 6).Value = 'High'
$ws1.Cells.Item(6, 7).Value = 'SW'
$ws1.Cells.Item(6, 8).Value = 1.556762700000007
$ws1.Cells.Item(6, 9).Value = 50.43
$ws1.Cells.Item(6, 10).Value = 54.38
$ws1.Cells.Item(6, 11).Value = 12
$ws1.Cells.Item(6, 12).Value = 1939
$ws1.Cells.Item(6, 13).Value = 'E'
$ws1.Cells.Item(6, 14).Value = 'E'
$ws1.Cells.Item(6, 15).Value = 31.94
$ws1.Cells.Item(6, 16).Value = 9.800000000000001
$ws1.Cells.Item(6, 17).Value = 'E'
$ws1.Cells.Item(6, 18).Value = 0.2
$ws1.Cells.Item(6, 19).Value = 0
$ws1.Cells.Item(6, 20).Value = 0
$ws1.Cells.Item(6, 21).Value = -2.2
$ws1.Cells.Item(6, 22).Value = '42.2860064, -85.6007573'
$ws1.Cells.Item(6, 23).Value = 43.5
$ws1.Cells.Item(6, 24).Value = -105
$ws1.Cells.Item(6, 25).Value = 43.5
$ws1.Cells.Item(6, 26).Value = -115
$ws1.Cells.Item(6, 31).Value = 0
$ws1.Cells.Item(6, 37).Value = '2025-12-03T05:15:42.459576'
$ws1.Cells.Item(6, 27).ClearContents()
$ws1.Cells.Item(6, 28).ClearContents()
$ws1.Cells.Item(6, 32).ClearContents()

# Row 7
$ws1.Cells.Item(7, 1).Value = 'Duke @ Virginia'
$ws1.Cells.Item(7, 2).Value = 'SAT 12/06'
$ws1.Cells.Item(7, 3).Value = '08:00 PM'
$ws1.Cells.Item(7, 4).Value = 'High'
$ws1.Cells.Item(7, 5).Value = 'NW-SE'
$ws1.Cells.Item(7, 6).Value = 'Med'
$ws1.Cells.Item(7, 8).Value = 67.46492769999999
$ws1.Cells.Item(7, 9).Value = 58.23
$ws1.Cells.Item(7, 10).Value = 61.08
$ws1.Cells.Item(7, 11).Value = 4.6
$ws1.Cells.Item(7, 12).Value = 1931
$ws1.Cells.Item(7, 13).Value = 'ENE'
$ws1.Cells.Item(7, 14).Value = 'E'
$ws1.Cells.Item(7, 15).Value = 39.2
$ws1.Cells.Item(7, 16).Value = 2.2
$ws1.Cells.Item(7, 17).Value = 'ENE'
$ws1.Cells.Item(7, 18).Value = 0
$ws1.Cells.Item(7, 19).Value = 0
$ws1.Cells.Item(7, 20).Value = 0
$ws1.Cells.Item(7, 21).Value = -2.4
$ws1.Cells.Item(7, 22).Value = '38.0311801, -78.5137897'
$ws1.Cells.Item(7, 23).Value = 58.5
$ws1.Cells.Item(7, 24).Value = -115
$ws1.Cells.Item(7, 25).Value = 57.5
$ws1.Cells.Item(7, 26).Value = -105
$ws1.Cells.Item(7, 27).Value = -3
$ws1.Cells.Item(7, 28).Value = -3.5
$ws1.Cells.Item(7, 31).Value = -0.0170940170940171
$ws1.Cells.Item(7, 32).Value = 0.5
$ws1.Cells.Item(7, 37).Value = '2025-12-03T05:15:42.459576'

# ---- Other sheet ----
# Row 2
$ws2.Cells.Item(2, 1).Value = 'North Dakota vs Tarleton State'
$ws2.Cells.Item(2, 2).Value = 'Tarleton State'
$ws2.Cells.Item(2, 3).Value = 'North Dakota'
$ws2.Cells.Item(2, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(2, 5).Value = '12:00 PM'
$ws2.Cells.Item(2, 6).Value = 'Low'
$ws2.Cells.Item(2, 11).Value = 66.93000000000001
$ws2.Cells.Item(2, 12).Value = 42.5
$ws2.Cells.Item(2, 15).Value = 'ENE'
$ws2.Cells.Item(2, 16).Value = 'ENE'
$ws2.Cells.Item(2, 17).Value = 68.78000000000002
$ws2.Cells.Item(2, 18).Value = 20.9
$ws2.Cells.Item(2, 19).Value = 'ENE'
$ws2.Cells.Item(2, 20).Value = 0
$ws2.Cells.Item(2, 21).Value = -6.5
$ws2.Cells.Item(2, 22).Value = 0
$ws2.Cells.Item(2, 24).Value = '32.2191836, -98.2130634'

# Row 3
$ws2.Cells.Item(3, 1).Value = 'South Dakota vs Mercer'
$ws2.Cells.Item(3, 2).Value = 'Mercer'
$ws2.Cells.Item(3, 3).Value = 'South Dakota'
$ws2.Cells.Item(3, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(3, 5).Value = '12:00 PM'
$ws2.Cells.Item(3, 6).Value = 'Low'
$ws2.Cells.Item(3, 10).Value = -234.7229156
$ws2.Cells.Item(3, 11).Value = 64.83
$ws2.Cells.Item(3, 12).Value = 51.08
$ws2.Cells.Item(3, 14).Value = 2013
$ws2.Cells.Item(3, 15).Value = 'NNW'
$ws2.Cells.Item(3, 16).Value = 'S'
$ws2.Cells.Item(3, 17).Value = 51.44
$ws2.Cells.Item(3, 18).Value = 2.2
$ws2.Cells.Item(3, 19).Value = 'S'
$ws2.Cells.Item(3, 20).Value = 2.5
$ws2.Cells.Item(3, 21).Value = -1.5
$ws2.Cells.Item(3, 22).Value = 0
$ws2.Cells.Item(3, 24).Value = '32.8262075, -83.6522485'

# Row 4
$ws2.Cells.Item(4, 1).Value = 'South Dakota State vs Montana'
$ws2.Cells.Item(4, 2).Value = 'Montana'
$ws2.Cells.Item(4, 3).Value = 'South Dakota State'
$ws2.Cells.Item(4, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(4, 5).Value = '12:00 PM'
$ws2.Cells.Item(4, 6).Value = 'High'
$ws2.Cells.Item(4, 10).Value = 474.5684815
$ws2.Cells.Item(4, 11).Value = 47.64
$ws2.Cells.Item(4, 12).Value = 46.7
$ws2.Cells.Item(4, 14).Value = 1986
$ws2.Cells.Item(4, 15).Value = 'E'
$ws2.Cells.Item(4, 16).Value = 'ESE'
$ws2.Cells.Item(4, 17).Value = 36.13999999999999
$ws2.Cells.Item(4, 18).Value = 2.5
$ws2.Cells.Item(4, 19).Value = 'ESE'
$ws2.Cells.Item(4, 20).Value = 1.2
$ws2.Cells.Item(4, 21).Value = -1.5
$ws2.Cells.Item(4, 22).Value = 0
$ws2.Cells.Item(4, 24).Value = '46.8638753, -113.9815042'

# Row 5
$ws2.Cells.Item(5, 1).Value = 'Villanova vs Lehigh'
$ws2.Cells.Item(5, 2).Value = 'Lehigh'
$ws2.Cells.Item(5, 3).Value = 'Villanova'
$ws2.Cells.Item(5, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(5, 5).Value = '12:00 PM'
$ws2.Cells.Item(5, 6).Value = 'High'
$ws2.Cells.Item(5, 10).Value = -37.06062315000001
$ws2.Cells.Item(5, 11).Value = 54.29
$ws2.Cells.Item(5, 12).Value = 55.05
$ws2.Cells.Item(5, 14).Value = 1988
$ws2.Cells.Item(5, 15).Value = 'ENE'
$ws2.Cells.Item(5, 16).Value = 'ENE'
$ws2.Cells.Item(5, 17).Value = 33.91999999999999
$ws2.Cells.Item(5, 18).Value = 5
$ws2.Cells.Item(5, 19).Value = 'ENE'
$ws2.Cells.Item(5, 20).Value = 0
$ws2.Cells.Item(5, 21).Value = 0
$ws2.Cells.Item(5, 22).Value = 0
$ws2.Cells.Item(5, 24).Value = '40.5890837, -75.3553874'

# Row 6
$ws2.Cells.Item(6, 1).Value = 'Abilene Christian vs Stephen F. Austin'
$ws2.Cells.Item(6, 2).Value = 'Stephen F. Austin'
$ws2.Cells.Item(6, 3).Value = 'Abilene Christian'
$ws2.Cells.Item(6, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(6, 5).Value = '01:00 PM'
$ws2.Cells.Item(6, 6).Value = 'Low'
$ws2.Cells.Item(6, 10).Value = -430.90566101
$ws2.Cells.Item(6, 11).Value = 68.06999999999999
$ws2.Cells.Item(6, 12).Value = 67.58
$ws2.Cells.Item(6, 14).Value = 1973
$ws2.Cells.Item(6, 15).Value = 'N'
$ws2.Cells.Item(6, 16).Value = 'N'
$ws2.Cells.Item(6, 17).Value = 63.14000000000001
$ws2.Cells.Item(6, 18).Value = 5.9
$ws2.Cells.Item(6, 19).Value = 'N'
$ws2.Cells.Item(6, 20).Value = 0
$ws2.Cells.Item(6, 21).Value = 0
$ws2.Cells.Item(6, 22).Value = 0
$ws2.Cells.Item(6, 24).Value = '31.625719, -94.6444034'

# Row 7
$ws2.Cells.Item(7, 1).Value = 'Yale vs Montana State'
$ws2.Cells.Item(7, 2).Value = 'Montana State'
$ws2.Cells.Item(7, 3).Value = 'Yale'
$ws2.Cells.Item(7, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(7, 5).Value = '12:00 PM'
$ws2.Cells.Item(7, 6).Value = 'High'
$ws2.Cells.Item(7, 10).Value = 1502.206045159
$ws2.Cells.Item(7, 11).Value = 42.68
$ws2.Cells.Item(7, 12).Value = 53.64
$ws2.Cells.Item(7, 14).Value = 1973
$ws2.Cells.Item(7, 15).Value = 'SW'
$ws2.Cells.Item(7, 16).Value = 'S'
$ws2.Cells.Item(7, 17).Value = 30.68
$ws2.Cells.Item(7, 18).Value = 3.1
$ws2.Cells.Item(7, 19).Value = 'S'
$ws2.Cells.Item(7, 20).Value = 1
$ws2.Cells.Item(7, 21).Value = 0
$ws2.Cells.Item(7, 22).Value = -3.5
$ws2.Cells.Item(7, 24).Value = '45.659048, -111.049547'

# Row 8
$ws2.Cells.Item(8, 1).Value = 'Rhode Island vs UC Davis'
$ws2.Cells.Item(8, 2).Value = 'UC Davis'
$ws2.Cells.Item(8, 3).Value = 'Rhode Island'
$ws2.Cells.Item(8, 4).Value = 'SAT 12/06'
$ws2.Cells.Item(8, 5).Value = '07:00 PM'
$ws2.Cells.Item(8, 6).Value = 'High'
$ws2.Cells.Item(8, 10).Value = -21.30542278
$ws2.Cells.Item(8, 11).Value = 62.21
$ws2.Cells.Item(8, 12).Value = 52.81
$ws2.Cells.Item(8, 14).Value = 2007
$ws2.Cells.Item(8, 15).Value = 'ESE'
$ws2.Cells.Item(8, 16).Value = 'SSE'
$ws2.Cells.Item(8, 17).Value = 54.5
$ws2.Cells.Item(8, 18).Value = 2.6
$ws2.Cells.Item(8, 19).Value = 'SE'
$ws2.Cells.Item(8, 20).Value = 0
$ws2.Cells.Item(8, 21).Value = 0
$ws2.Cells.Item(8, 22).Value = 0
$ws2.Cells.Item(8, 24).Value = '38.5365266, -121.7627936'

